$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header: "Dev2.Integration.Tests" (bold, matches style of existing
# section headers in A1 / A14).
$ws.Range("A18").Value = "Dev2.Integration.Tests"
$ws.Range("A18").Font.Bold = $true

# Row 19: Add_A_New_InputOnSqlProcedure_Expect_New_IS_InputAdded
$ws.Range("A19").Value = "Add_A_New_InputOnSqlProcedure_Expect_New_IS_InputAdded"
$ws.Range("B19").Value = "copy resources, restart server, run again"
$ws.Range("B19").WrapText = $true

# Row 20: Change_sql_source_verify_Empty_Inputs
$ws.Range("A20").Value = "Change_sql_source_verify_Empty_Inputs"
$ws.Range("B20").Value = "copy resources, restart server, run again"
$ws.Range("B20").WrapText = $true

# Row 21: ExecutionWithNoStartNode_ExpectedInvalidValidResult
$ws.Range("A21").Value = "ExecutionWithNoStartNode_ExpectedInvalidValidResult"
$ws.Range("B21").Value = "copy file WorkflowWithNoStartNodeConnected.xml in Resources\Acceptance Testing Resources directory, restart server , run again"
$ws.Range("B21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 28.8

# Update view state: scroll so row 7 is at the top and select B22 (the cell
# just below the newly added rows), matching where the author's cursor ended
# up after entering the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B22").Select()
